$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3 (PORCELANATO)
$ws.Range("D3").Value = 809.04
$ws.Range("E3").Value = 12914.3
$ws.Range("F3").Value = 0.05895357835628935

# Row 4 (TOTAL)
$ws.Range("D4").Value = 2864.98
$ws.Range("E4").Value = 10858.36
$ws.Range("F4").Value = 0.2087669619786437
